# Closes #5820 - uniform heading colors + resized/restyled heading levels.
$d = $word.ActiveDocument

# Heading 1: drop the 35% theme-shade on the accent1 color so it matches
# the other headings' plain accent1 color (4F81BD).
$h1 = $d.Styles("Heading 1")
$h1.Font.Color = -738131969   # wdColor sentinel: accent1 theme color, no shade

# Heading 2: 16pt -> 14pt
$h2 = $d.Styles("Heading 2")
$h2.Font.Size = 14
$h2.Font.SizeBi = 14

# Heading 3: 14pt -> 12pt
$h3 = $d.Styles("Heading 3")
$h3.Font.Size = 12
$h3.Font.SizeBi = 12

# Heading 4: bold -> italic
$h4 = $d.Styles("Heading 4")
$h4.Font.Bold = $false
$h4.Font.Italic = $true

# Heading 5: no longer italic
$h5 = $d.Styles("Heading 5")
$h5.Font.Italic = $false
